# Update scripts to output nr of studies and effect sizes for each moderator level.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Update "descriptives" sheet (sheet1) with refreshed model results
# ---------------------------------------------------------------------------
$descriptives = $wb.Worksheets.Item("descriptives")

# Row 2 - NS
$descriptives.Range("C2").Value = 845
$descriptives.Range("D2").Value = 0.008498832191968928
$descriptives.Range("E2").Value = 0.05055652050667844
$descriptives.Range("F2").Value = 97.74939134949946
$descriptives.Range("G2").Value = 83.68198448486274
$descriptives.Range("H2").Value = 14.06740686463673

# Row 3 - NT
$descriptives.Range("C3").Value = 438
$descriptives.Range("D3").Value = 0.003358330498989658
$descriptives.Range("E3").Value = 0.05303031348248766
$descriptives.Range("F3").Value = 95.05703761436311
$descriptives.Range("G3").Value = 89.39573906161222
$descriptives.Range("H3").Value = 5.66129855275089

# ---------------------------------------------------------------------------
# 2) Update "coefficients" sheet (sheet2) with refreshed model results
# ---------------------------------------------------------------------------
$coefficients = $wb.Worksheets.Item("coefficients")

# Row 2 - NS / intercept
$coefficients.Range("D2").Value = 0.04639171878434833
$coefficients.Range("E2").Value = 0.06660707138845227
$coefficients.Range("F2").Value = 0.6969987142224494
$coefficients.Range("G2").Value = 0.491887878560833
$coefficients.Range("H2").Value = -0.09013262919490939
$coefficients.Range("I2").Value = 0.181204681005399
$coefficients.Range("J2").Value = 26.43614495246367

# Row 3 - NS / quality_score_out_of_5
$coefficients.Range("D3").Value = 0.01755648388327493
$coefficients.Range("E3").Value = 0.01835977092344859
$coefficients.Range("F3").Value = 0.9563457029266192
$coefficients.Range("G3").Value = 0.3468955748784063
$coefficients.Range("H3").Value = -0.02000834332856446
$coefficients.Range("I3").Value = 0.05507181264810409
$coefficients.Range("J3").Value = 28.66008235053942

# Row 4 - NT / intercept
$coefficients.Range("D4").Value = 0.000504142035388867
$coefficients.Range("E4").Value = 0.07945720451557985
$coefficients.Range("F4").Value = 0.006344825257485614
$coefficients.Range("G4").Value = 0.9950089843215419
$coefficients.Range("H4").Value = -0.1651784088797527
$coefficients.Range("I4").Value = 0.1661590194552117
$coefficients.Range("J4").Value = 17.5941029536253

# Row 5 - NT / quality_score_out_of_5
$coefficients.Range("D5").Value = 0.0127506238325392
$coefficients.Range("E5").Value = 0.02371268760343636
$coefficients.Range("F5").Value = 0.5377422882316016
$coefficients.Range("G5").Value = 0.5969555672156028
$coefficients.Range("H5").Value = -0.03684020372310787
$coefficients.Range("I5").Value = 0.06227880647813155
$coefficients.Range("J5").Value = 19.13119962178465

# ---------------------------------------------------------------------------
# 3) Add new "nr_studies" sheet (sheet3) at the end of the workbook
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$nrStudies = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$nrStudies.Name = "nr_studies"

# Header row
$nrStudies.Cells.Item(1, 1).Value = "outcome"
$nrStudies.Cells.Item(1, 2).Value = "quality_score_out_of_5"
$nrStudies.Cells.Item(1, 3).Value = "n_effect_sizes"
$nrStudies.Cells.Item(1, 4).Value = "k_studies"

$headerRange = $nrStudies.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108  # xlCenter

# Data rows
$data = @(
    @("NS", 3, 300, 31),
    @("NS", 4, 393, 31),
    @("NS", 2, 74, 8),
    @("NS", 5, 69, 9),
    @("NS", 1, 9, 3),
    @("NT", 3, 149, 18),
    @("NT", 4, 195, 21),
    @("NT", 2, 41, 4),
    @("NT", 5, 53, 8)
)

$rowIndex = 2
foreach ($entry in $data) {
    $nrStudies.Cells.Item($rowIndex, 1).Value = $entry[0]
    $nrStudies.Cells.Item($rowIndex, 2).Value = $entry[1]
    $nrStudies.Cells.Item($rowIndex, 3).Value = $entry[2]
    $nrStudies.Cells.Item($rowIndex, 4).Value = $entry[3]
    $rowIndex++
}

$descriptives.Select()
